$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amplicon")
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
Write-Host ("After add: " + $lo.Range.Address())
foreach ($col in $lo.ListColumns) {
    Write-Host ($col.Index.ToString() + " " + $col.Name)
}
